$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for the columns that change: D (Fecha), J (Volumen),
# K (Precio mínimo), L (Precio máximo), M (Precio promedio ponderado), P (Precio $/Kg)
# for rows 2..5, then perform a cyclic shift: row2->row3, row3->row4, row4->row5, row5->row2.

$rows = 2..5
$cols = @("D", "J", "K", "L", "M", "P")

$original = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $original[$r] = $rowVals
}

foreach ($r in $rows) {
    # source row is the one above (wrapping from 2 to 5)
    $src = $r - 1
    if ($src -lt 2) { $src = 5 }
    $rowVals = $original[$src]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}
